$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Cells.Item(2, 8).Value = 1017.05554
$ws.Cells.Item(2, 9).Value = 852.5455
$ws.Cells.Item(2, 11).Value = 852.5455
$ws.Cells.Item(2, 13).Value = -739.5455

# Row 18 (Leve Item ID 5471)
$ws.Cells.Item(18, 8).Value = 1258.0
$ws.Cells.Item(18, 9).Value = 1258.0
$ws.Cells.Item(18, 11).Value = 1258.0
$ws.Cells.Item(18, 13).Value = -974.0

# Row 33 (Leve Item ID 5512)
$ws.Cells.Item(33, 8).Value = 639889.94
$ws.Cells.Item(33, 9).Value = 862591.56
$ws.Cells.Item(33, 10).Value = 3599.7144
$ws.Cells.Item(33, 11).Value = 862591.56
$ws.Cells.Item(33, 12).Value = 3599.7144
$ws.Cells.Item(33, 13).Value = -862362.56
$ws.Cells.Item(33, 14).Value = -4057.7144

# Row 42 (Leve Item ID 4600)
$ws.Cells.Item(42, 8).Value = 615.05884
$ws.Cells.Item(42, 9).Value = 76.666664
$ws.Cells.Item(42, 10).Value = 908.7273
$ws.Cells.Item(42, 11).Value = 229.999992
$ws.Cells.Item(42, 12).Value = 2726.1819
$ws.Cells.Item(42, 13).Value = 0.000008000000008223651
$ws.Cells.Item(42, 14).Value = -3186.1819

# Row 70 (Leve Item ID 12604)
$ws.Cells.Item(70, 8).Value = 12277.777
$ws.Cells.Item(70, 9).Value = 12277.777
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 36833.331
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = -36563.331
$ws.Cells.Item(70, 14).ClearContents()

# Row 73 (Leve Item ID 12604)
$ws.Cells.Item(73, 8).Value = 12277.777
$ws.Cells.Item(73, 9).Value = 12277.777
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 36833.331
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = -35897.331
$ws.Cells.Item(73, 14).ClearContents()

# Row 74 (Leve Item ID 5507)
$ws.Cells.Item(74, 8).Value = 4118.091
$ws.Cells.Item(74, 9).Value = 2259.8
$ws.Cells.Item(74, 10).Value = 5666.6665
$ws.Cells.Item(74, 11).Value = 2259.8
$ws.Cells.Item(74, 12).Value = 5666.6665
$ws.Cells.Item(74, 13).Value = -1323.8
$ws.Cells.Item(74, 14).Value = -7538.6665

# Row 77 (Leve Item ID 5507)
$ws.Cells.Item(77, 8).Value = 4118.091
$ws.Cells.Item(77, 9).Value = 2259.8
$ws.Cells.Item(77, 10).Value = 5666.6665
$ws.Cells.Item(77, 11).Value = 11299.0
$ws.Cells.Item(77, 12).Value = 28333.3325
$ws.Cells.Item(77, 13).Value = -6619.0
$ws.Cells.Item(77, 14).Value = -37693.3325

# Row 106 (Leve Item ID 19903)
$ws.Cells.Item(106, 8).Value = 10291.923
$ws.Cells.Item(106, 9).Value = 1538.0
$ws.Cells.Item(106, 11).Value = 1538.0
$ws.Cells.Item(106, 13).Value = -907.0

# Row 127 (Leve Item ID 36114)
$ws.Cells.Item(127, 8).Value = 52325.39
$ws.Cells.Item(127, 9).Value = 52325.39
$ws.Cells.Item(127, 11).Value = 156976.17
$ws.Cells.Item(127, 13).Value = -152016.17

# Row 135 (Leve Item ID 44047)
$ws.Cells.Item(135, 8).Value = 691.64703
$ws.Cells.Item(135, 9).Value = 648.6667
$ws.Cells.Item(135, 10).Value = 1014.0
$ws.Cells.Item(135, 11).Value = 5838.0003
$ws.Cells.Item(135, 12).Value = 9126.0
$ws.Cells.Item(135, 13).Value = -3303.0003
$ws.Cells.Item(135, 14).Value = -14196.0

# Row 138 (Leve Item ID 44169)
$ws.Cells.Item(138, 8).Value = 3262.9583
$ws.Cells.Item(138, 9).Value = 2244.8333
$ws.Cells.Item(138, 10).Value = 3602.3333
$ws.Cells.Item(138, 11).Value = 6734.499899999999
$ws.Cells.Item(138, 12).Value = 10806.9999
$ws.Cells.Item(138, 13).Value = -1594.499899999999
$ws.Cells.Item(138, 14).Value = -21086.9999

# Row 141 (Leve Item ID 44161)
$ws.Cells.Item(141, 8).Value = 42603.727
$ws.Cells.Item(141, 9).Value = 42603.727
$ws.Cells.Item(141, 11).Value = 127811.181
$ws.Cells.Item(141, 13).Value = -122631.181

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 12218482.0
$ws.Cells.Item(32, 9).Value = 2462448.0
$ws.Cells.Item(32, 10).Value = 37150570.0
$ws.Cells.Item(32, 11).Value = 2462448.0
$ws.Cells.Item(32, 12).Value = 37150570.0
$ws.Cells.Item(32, 13).Value = -2462161.0
$ws.Cells.Item(32, 14).Value = -37151144.0

# Row 110 (Leve Item ID 27708)
$ws.Cells.Item(110, 8).Value = 1020.2
$ws.Cells.Item(110, 9).Value = 768.64514
$ws.Cells.Item(110, 10).Value = 2969.75
$ws.Cells.Item(110, 11).Value = 768.64514
$ws.Cells.Item(110, 12).Value = 2969.75
$ws.Cells.Item(110, 13).Value = 1276.35486
$ws.Cells.Item(110, 14).Value = -7059.75

# Row 122 (Leve Item ID 36168)
$ws.Cells.Item(122, 8).Value = 26047.223
$ws.Cells.Item(122, 9).Value = 36604.168
$ws.Cells.Item(122, 11).Value = 109812.504
$ws.Cells.Item(122, 13).Value = -107362.504

# Row 132 (Leve Item ID 43997)
$ws.Cells.Item(132, 8).Value = 2659.4707
$ws.Cells.Item(132, 9).Value = 2498.4773
$ws.Cells.Item(132, 11).Value = 7495.4319
$ws.Cells.Item(132, 13).Value = -4965.4319

$ws = $wb.Worksheets.Item("BSM")
# Row 26 (Leve Item ID 19535)
$ws.Cells.Item(26, 8).Value = 15339.667
$ws.Cells.Item(26, 9).Value = 7750.0
$ws.Cells.Item(26, 10).Value = 30519.0
$ws.Cells.Item(26, 11).Value = 7750.0
$ws.Cells.Item(26, 12).Value = 30519.0
$ws.Cells.Item(26, 13).Value = -7458.0
$ws.Cells.Item(26, 14).Value = -31103.0

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (Leve Item ID 3742)
$ws.Cells.Item(4, 8).Value = 4384746.5
$ws.Cells.Item(4, 9).Value = 11138.786
$ws.Cells.Item(4, 10).Value = 35000000.0
$ws.Cells.Item(4, 11).Value = 11138.786
$ws.Cells.Item(4, 12).Value = 35000000.0
$ws.Cells.Item(4, 13).Value = -11026.786
$ws.Cells.Item(4, 14).Value = -35000224.0

# Row 31 (Leve Item ID 44023)
$ws.Cells.Item(31, 8).Value = 3408.1072
$ws.Cells.Item(31, 9).Value = 1628.9286
$ws.Cells.Item(31, 11).Value = 1628.9286
$ws.Cells.Item(31, 13).Value = -1333.9286

# Row 33 (Leve Item ID 1836)
$ws.Cells.Item(33, 8).Value = 1980.0
$ws.Cells.Item(33, 9).Value = 1650.0
$ws.Cells.Item(33, 10).Value = 2200.0
$ws.Cells.Item(33, 11).Value = 1650.0
$ws.Cells.Item(33, 12).Value = 2200.0
$ws.Cells.Item(33, 13).Value = -1271.0
$ws.Cells.Item(33, 14).Value = -2958.0

# Row 34 (Leve Item ID 44023)
$ws.Cells.Item(34, 8).Value = 3408.1072
$ws.Cells.Item(34, 9).Value = 1628.9286
$ws.Cells.Item(34, 11).Value = 1628.9286
$ws.Cells.Item(34, 13).Value = -1426.9286

# Row 58 (Leve Item ID 44021)
$ws.Cells.Item(58, 8).Value = 1962.25
$ws.Cells.Item(58, 9).Value = 2400.0
$ws.Cells.Item(58, 10).Value = 649.0
$ws.Cells.Item(58, 11).Value = 2400.0
$ws.Cells.Item(58, 12).Value = 649.0
$ws.Cells.Item(58, 13).Value = -2197.0
$ws.Cells.Item(58, 14).Value = -1055.0

# Row 62 (Leve Item ID 12580)
$ws.Cells.Item(62, 8).Value = 5438.952
$ws.Cells.Item(62, 10).Value = 5979.4443
$ws.Cells.Item(62, 12).Value = 5979.4443
$ws.Cells.Item(62, 14).Value = -7227.4443

# Row 65 (Leve Item ID 12580)
$ws.Cells.Item(65, 8).Value = 5438.952
$ws.Cells.Item(65, 10).Value = 5979.4443
$ws.Cells.Item(65, 12).Value = 29897.2215
$ws.Cells.Item(65, 14).Value = -36137.2215

# Row 97 (Leve Item ID 19730)
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()

# Row 132 (Leve Item ID 44019)
$ws.Cells.Item(132, 8).Value = 3450.6177
$ws.Cells.Item(132, 9).Value = 3316.5
$ws.Cells.Item(132, 11).Value = 9949.5
$ws.Cells.Item(132, 13).Value = -7419.5

# Row 136 (Leve Item ID 44021)
$ws.Cells.Item(136, 8).Value = 1962.25
$ws.Cells.Item(136, 9).Value = 2400.0
$ws.Cells.Item(136, 10).Value = 649.0
$ws.Cells.Item(136, 11).Value = 7200.0
$ws.Cells.Item(136, 12).Value = 1947.0
$ws.Cells.Item(136, 13).Value = -4650.0
$ws.Cells.Item(136, 14).Value = -7047.0

$ws = $wb.Worksheets.Item("CUL")
# Row 12 (Leve Item ID 4854)
$ws.Cells.Item(12, 8).Value = 1000.0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 1000.0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 3000.0
$ws.Cells.Item(12, 13).ClearContents()
$ws.Cells.Item(12, 14).Value = -3346.0

# Row 102 (Leve Item ID 19813)
$ws.Cells.Item(102, 8).Value = 9653.23
$ws.Cells.Item(102, 10).Value = 10041.0
$ws.Cells.Item(102, 12).Value = 30123.0
$ws.Cells.Item(102, 14).Value = -34991.0

# Row 112 (Leve Item ID 27855)
$ws.Cells.Item(112, 8).Value = 201851.2
$ws.Cells.Item(112, 9).Value = 251306.5
$ws.Cells.Item(112, 11).Value = 753919.5
$ws.Cells.Item(112, 13).Value = -752811.5

# Row 113 (Leve Item ID 27843)
$ws.Cells.Item(113, 8).Value = 828.0714
$ws.Cells.Item(113, 9).Value = 771.0
$ws.Cells.Item(113, 10).Value = 904.1667
$ws.Cells.Item(113, 11).Value = 2313.0
$ws.Cells.Item(113, 12).Value = 2712.5001
$ws.Cells.Item(113, 13).Value = -143.0
$ws.Cells.Item(113, 14).Value = -7052.5001

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Cells.Item(80, 8).Value = 8075.3335
$ws.Cells.Item(80, 9).Value = 12047.7
$ws.Cells.Item(80, 10).Value = 3109.875
$ws.Cells.Item(80, 11).Value = 12047.7
$ws.Cells.Item(80, 12).Value = 3109.875
$ws.Cells.Item(80, 13).Value = -11049.7
$ws.Cells.Item(80, 14).Value = -5105.875

# Row 83 (Leve Item ID 12521)
$ws.Cells.Item(83, 8).Value = 8075.3335
$ws.Cells.Item(83, 9).Value = 12047.7
$ws.Cells.Item(83, 10).Value = 3109.875
$ws.Cells.Item(83, 11).Value = 60238.5
$ws.Cells.Item(83, 12).Value = 15549.375
$ws.Cells.Item(83, 13).Value = -55246.5
$ws.Cells.Item(83, 14).Value = -25533.375

# Row 126 (Leve Item ID 36184)
$ws.Cells.Item(126, 8).Value = 323142.0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 13).ClearContents()

# Row 132 (Leve Item ID 44008)
$ws.Cells.Item(132, 8).Value = 3697.5
$ws.Cells.Item(132, 9).Value = 4030.1155
$ws.Cells.Item(132, 10).Value = 3157.0
$ws.Cells.Item(132, 11).Value = 12090.3465
$ws.Cells.Item(132, 12).Value = 9471.0
$ws.Cells.Item(132, 13).Value = -9560.3465
$ws.Cells.Item(132, 14).Value = -14531.0

$ws = $wb.Worksheets.Item("LTW")
# Row 100 (Leve Item ID 19995)
$ws.Cells.Item(100, 8).Value = 29249.621
$ws.Cells.Item(100, 9).Value = 21618.65
$ws.Cells.Item(100, 11).Value = 21618.65
$ws.Cells.Item(100, 13).Value = -21077.65

# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 3703.6667
$ws.Cells.Item(132, 9).Value = 3400.16
$ws.Cells.Item(132, 11).Value = 10200.48
$ws.Cells.Item(132, 13).Value = -7670.48

# Row 136 (Leve Item ID 44060)
$ws.Cells.Item(136, 8).Value = 3484.0833
$ws.Cells.Item(136, 9).Value = 3003.4285
$ws.Cells.Item(136, 10).Value = 4157.0
$ws.Cells.Item(136, 11).Value = 9010.2855
$ws.Cells.Item(136, 12).Value = 12471.0
$ws.Cells.Item(136, 13).Value = -6460.2855
$ws.Cells.Item(136, 14).Value = -17571.0

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (Leve Item ID 12596)
$ws.Cells.Item(81, 8).Value = 51262.117
$ws.Cells.Item(81, 9).Value = 985.8461
$ws.Cells.Item(81, 11).Value = 1971.6922
$ws.Cells.Item(81, 13).Value = -910.6922

# Row 84 (Leve Item ID 12596)
$ws.Cells.Item(84, 8).Value = 51262.117
$ws.Cells.Item(84, 9).Value = 985.8461
$ws.Cells.Item(84, 11).Value = 9858.461
$ws.Cells.Item(84, 13).Value = -4554.460999999999

# Row 100 (Leve Item ID 19981)
$ws.Cells.Item(100, 8).Value = 2534.1765
$ws.Cells.Item(100, 9).Value = 3488.1
$ws.Cells.Item(100, 11).Value = 6976.2
$ws.Cells.Item(100, 13).Value = -6435.2

# Row 132 (Leve Item ID 44029)
$ws.Cells.Item(132, 8).Value = 8198.474
$ws.Cells.Item(132, 9).Value = 8457.588
$ws.Cells.Item(132, 10).Value = 5996.0
$ws.Cells.Item(132, 11).Value = 25372.764
$ws.Cells.Item(132, 12).Value = 17988.0
$ws.Cells.Item(132, 13).Value = -22842.764
$ws.Cells.Item(132, 14).Value = -23048.0

# Row 136 (Leve Item ID 44031)
$ws.Cells.Item(136, 8).Value = 8351.1
$ws.Cells.Item(136, 9).Value = 5148.2144
$ws.Cells.Item(136, 10).Value = 15824.5
$ws.Cells.Item(136, 11).Value = 15444.6432
$ws.Cells.Item(136, 12).Value = 47473.5
$ws.Cells.Item(136, 13).Value = -12894.6432
$ws.Cells.Item(136, 14).Value = -52573.5
